$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("S5 Table")

$data = @(
    @(0.64951644100580297, 0.51726155646576899, 0.487910189982729),
    @(0.64676034348165501, 0.51639344262295095, 0.48700173310225298),
    @(0.66908563134978205, 0.53667481662591698, 0.44153005464480899),
    @(0.68548039638086999, 0.553968253968254, 0.434223541048467),
    @(0.70083333333333298, 0.56008251676121701, 0.434108527131783),
    @(0.68105158730158699, 0.53891050583657596, 0.49306296691568802),
    @(0.86853817185895199, 0.79839357429718905, 0.55908639523336701),
    @(0.87579821200510899, 0.80769230769230804, 0.55908639523336701),
    @(0.86827859978347199, 0.79609634551495001, 0.55908639523336701),
    @(0.87081339712918704, 0.79478553406223695, 0.55908639523336701),
    @(0.86853817185895199, 0.79839357429718905, 0.55908639523336701),
    @(0.87875288683602804, 0.73941987636709505, 0.55129434324065196),
    @(0.87875288683602804, 0.73941987636709505, 0.55129434324065196),
    @(0.88184998056743102, 0.75071225071225101, 0.55404178019981798),
    @(0.88169123351435197, 0.75047438330170801, 0.55525362318840599),
    @(0.88399528116397996, 0.74709302325581395, 0.55636363636363595),
    @(0.88569293985804998, 0.75918367346938798, 0.55212677231025897),
    @(0.88556469708302199, 0.75840145322434205, 0.55379482902418697),
    @(0.88413479742521806, 0.753361149745016, 0.55029093931837103),
    @(0.87070091423595997, 0.78942353501667495, 0.58310991957104596),
    @(0.87177280550774505, 0.79194001874414199, 0.58289817232375996),
    @(0.87187362878455499, 0.78884078884078901, 0.58155540261527905),
    @(0.86734693877550995, 0.77528675703858196, 0.57097032878909404),
    @(0.86903729401561203, 0.78802281368821303, 0.58404746209624303),
    @(0.87494736842105303, 0.72616525423728795, 0.51903807615230502),
    @(0.90670220868240703, 0.78614982578397197, 0.66754155730533704),
    @(0.90858623242042902, 0.78264568268821799, 0.660519698239732)
)

$startRow = 4
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value = $data[$i][0]
    $ws.Cells.Item($r, 3).Value = $data[$i][1]
    $ws.Cells.Item($r, 4).Value = $data[$i][2]
}

$ws.Range("B4:B30").Select()
